$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 (columns E, G, H, K, L, M, N, O, P, Q, R, S, T)

# Row 2 (ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3884013333333334
$ws.Range("H2").Value = 1.165204
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.139245
$ws.Range("N2").Value = 0.417735
$ws.Range("O2").Value = 0.01212793695325064
$ws.Range("P2").Value = 0.01283499108585158
$ws.Range("Q2").Value = 0.05408294366
$ws.Range("R2").Value = 0.48674649294
$ws.Range("S2").Value = 0.01212793695325064
$ws.Range("T2").Value = 0.01283499108585158

# Row 3 (FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3884013333333334
$ws.Range("H3").Value = 1.165204
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.380691666666667
$ws.Range("N3").Value = 28.142075
$ws.Range("O3").Value = 0.8170378621222814
$ws.Range("P3").Value = 0.8646708601442703
$ws.Range("Q3").Value = 3.643473150922223
$ws.Range("R3").Value = 32.79125835830001
$ws.Range("S3").Value = 0.8170378621222814
$ws.Range("T3").Value = 0.8646708601442703

# Row 4 (now M1 -- target clusters were re-sorted, M1/M2 inserted before sCs)
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3884013333333334
$ws.Range("H4").Value = 1.165204
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.006356333333333333
$ws.Range("N4").Value = 0.019069
$ws.Range("O4").Value = 0.0005536228225107701
$ws.Range("P4").Value = 0.0005858988234553099
$ws.Range("Q4").Value = 0.002468808341777778
$ws.Range("R4").Value = 0.022219275076
$ws.Range("S4").Value = 0.0005536228225107701
$ws.Range("T4").Value = 0.0005858988234553099

# Row 5 (new, M2)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3884013333333334
$ws.Range("H5").Value = 1.165204
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05759666666666666
$ws.Range("N5").Value = 0.17279
$ws.Range("O5").Value = 0.005016544522609259
$ws.Range("P5").Value = 0.005309007168957103
$ws.Range("Q5").Value = 0.02237062212888889
$ws.Range("R5").Value = 0.20133559916
$ws.Range("S5").Value = 0.005016544522609259
$ws.Range("T5").Value = 0.005309007168957103

# Row 6 (new, sCs -- re-added after the M1/M2 split)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3884013333333334
$ws.Range("H6").Value = 1.165204
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.897453
$ws.Range("N6").Value = 3.794906
$ws.Range("O6").Value = 0.1652640335793479
$ws.Range("P6").Value = 0.1165992427774658
$ws.Range("Q6").Value = 0.7369732751373335
$ws.Range("R6").Value = 4.421839650824
$ws.Range("S6").Value = 0.1652640335793479
$ws.Range("T6").Value = 0.1165992427774658
